# EU events exog and reduced cap
# Update inventory_cases_end_of_month (column E) values on the "Main" sheet
# for the February 2026 reporting rows, reflecting refreshed query results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$updates = @{
    4  = 63
    8  = 32
    10 = 99
    13 = 36
    28 = 124
    36 = 40
    41 = 70
    42 = 1
    43 = 14
    53 = 100
    57 = 56
    63 = 137
    70 = 197
    71 = 103
    73 = 34
    74 = 8
    77 = 102
}

foreach ($row in $updates.Keys) {
    $ws.Range("E" + $row).Value = $updates[$row]
}
